# Weekly data refresh: insert a new daily record at the top of the data
# block (row 171), pushing all subsequent rows down by one (171-207 -> 172-208).
# The new row re-uses the static/unchanged attributes of the record that used
# to sit at row 171 (market, region, category, variety, quality, unit,
# classification, Kg/units) and carries its own new date/volume/price/origin
# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before row 171; everything from 171..207 shifts to 172..208
$ws.Rows("171").Insert()

# Populate the newly inserted row 171 with the new weekly record.
$ws.Cells.Item(171, 1).Value  = 8                                    # A Mercado ID
$ws.Cells.Item(171, 2).Value  = 'Terminal La Palmera de La Serena'   # B Mercado
$ws.Cells.Item(171, 3).Value  = 'Coquimbo'                           # C Región
$ws.Cells.Item(171, 4).Value  = 44511                                # D Fecha
$ws.Cells.Item(171, 5).Value  = 4                                    # E Codreg
$ws.Cells.Item(171, 6).Value  = 100112032                            # F Categoría ID
$ws.Cells.Item(171, 7).Value  = 'Zapallo italiano'                   # G Categoría
$ws.Cells.Item(171, 8).Value  = 'Sin especificar'                    # H Variedad
$ws.Cells.Item(171, 9).Value  = 'Primera'                            # I Calidad
$ws.Cells.Item(171, 10).Value = 500                                  # J Volumen
$ws.Cells.Item(171, 11).Value = 10000                                # K Precio mínimo
$ws.Cells.Item(171, 12).Value = 10500                                # L Precio máximo
$ws.Cells.Item(171, 13).Value = 10250                                # M Precio promedio ponderado
$ws.Cells.Item(171, 14).Value = '$/caja 70 unidades'                 # N Unidad de comercialización
$ws.Cells.Item(171, 15).Value = 'Provincia de Limarí'                # O Origen
$ws.Cells.Item(171, 16).Value = 146                                  # P Precio $/Kg
$ws.Cells.Item(171, 17).Value = 70                                   # Q Kg o Unidades
$ws.Cells.Item(171, 18).Value = 'Hortaliza'                          # R Clasificación
